# Append 3 new match rows (123, 124, 125) to the Belgium Jupiler Pro League
# 2023-2024 sheet, mirroring the formatting of the existing last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (122) down onto
# the three new rows (123-125) before writing values, so the bold/bordered
# "Indice" style (column A) and the datetime number format (column E) carry
# over exactly as in the rest of the table.
$ws.Range("A122:V122").Copy()
$ws.Range("A123:V125").PasteSpecial(-4122)

# --- Row 123: RWDM 0 - 0 Charleroi ---
$ws.Cells.Item(123, 1).Value = 122
$ws.Cells.Item(123, 2).Value = "belgium"
$ws.Cells.Item(123, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(123, 4).Value = "2023-2024"
$ws.Cells.Item(123, 5).Value = 45262.66666666666
$ws.Cells.Item(123, 6).Value = "RWDM"
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = "Charleroi"
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 2.79
$ws.Cells.Item(123, 11).Value = "26/11/2023 18:43"
$ws.Cells.Item(123, 12).Value = 3.06
$ws.Cells.Item(123, 13).Value = "02/12/2023 15:57"
$ws.Cells.Item(123, 14).Value = 3.58
$ws.Cells.Item(123, 15).Value = "26/11/2023 18:43"
$ws.Cells.Item(123, 16).Value = 3.35
$ws.Cells.Item(123, 17).Value = "02/12/2023 15:57"
$ws.Cells.Item(123, 18).Value = 2.36
$ws.Cells.Item(123, 19).Value = "26/11/2023 18:43"
$ws.Cells.Item(123, 20).Value = 2.45
$ws.Cells.Item(123, 21).Value = "02/12/2023 15:57"
$ws.Cells.Item(123, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/rwd-molenbeek-charleroi/dWcYOUgf/"

# --- Row 124: Eupen 1 - 1 Kortrijk ---
$ws.Cells.Item(124, 1).Value = 123
$ws.Cells.Item(124, 2).Value = "belgium"
$ws.Cells.Item(124, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(124, 4).Value = "2023-2024"
$ws.Cells.Item(124, 5).Value = 45262.76041666666
$ws.Cells.Item(124, 6).Value = "Eupen"
$ws.Cells.Item(124, 7).Value = 1
$ws.Cells.Item(124, 8).Value = "Kortrijk"
$ws.Cells.Item(124, 9).Value = 1
$ws.Cells.Item(124, 10).Value = 2.11
$ws.Cells.Item(124, 11).Value = "26/11/2023 19:42"
$ws.Cells.Item(124, 12).Value = 2.42
$ws.Cells.Item(124, 13).Value = "02/12/2023 18:14"
$ws.Cells.Item(124, 14).Value = 3.81
$ws.Cells.Item(124, 15).Value = "26/11/2023 19:42"
$ws.Cells.Item(124, 16).Value = 3.46
$ws.Cells.Item(124, 17).Value = "02/12/2023 18:14"
$ws.Cells.Item(124, 18).Value = 3.07
$ws.Cells.Item(124, 19).Value = "26/11/2023 19:42"
$ws.Cells.Item(124, 20).Value = 3.03
$ws.Cells.Item(124, 21).Value = "02/12/2023 18:06"
$ws.Cells.Item(124, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/eupen-kortrijk/befkLngJ/"

# --- Row 125: Antwerp 1 - 0 Leuven ---
$ws.Cells.Item(125, 1).Value = 124
$ws.Cells.Item(125, 2).Value = "belgium"
$ws.Cells.Item(125, 3).Value = "jupiler-pro-league"
$ws.Cells.Item(125, 4).Value = "2023-2024"
$ws.Cells.Item(125, 5).Value = 45262.86458333334
$ws.Cells.Item(125, 6).Value = "Antwerp"
$ws.Cells.Item(125, 7).Value = 1
$ws.Cells.Item(125, 8).Value = "Leuven"
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 1.36
$ws.Cells.Item(125, 11).Value = "26/11/2023 16:13"
$ws.Cells.Item(125, 12).Value = 1.29
$ws.Cells.Item(125, 13).Value = "02/12/2023 20:44"
$ws.Cells.Item(125, 14).Value = 5.37
$ws.Cells.Item(125, 15).Value = "26/11/2023 16:13"
$ws.Cells.Item(125, 16).Value = 6.13
$ws.Cells.Item(125, 17).Value = "02/12/2023 20:44"
$ws.Cells.Item(125, 18).Value = 6.84
$ws.Cells.Item(125, 19).Value = "26/11/2023 16:13"
$ws.Cells.Item(125, 20).Value = 10.06
$ws.Cells.Item(125, 21).Value = "02/12/2023 20:44"
$ws.Cells.Item(125, 22).Value = "https://www.betexplorer.com/football/belgium/jupiler-pro-league/antwerp-leuven/00YXqWVJ/"
